# Update return shipments data (auto)
# Inserts a new "courier_name" column between carrier_slug (B) and
# status_tag (old C) and populates it with the courier display name
# for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C, shifting status_tag..custom_fields_json
# (old C:J) to D:K.
$ws.Columns.Item(3).Insert()

# New column header
$ws.Range("C1").Value = "courier_name"

# New column values (one per data row)
$ws.Range("C2").Value = "DB Schenker"
$ws.Range("C3").Value = "DHL Express"
$ws.Range("C4").Value = "kn"
$ws.Range("C5").Value = "testing-courier"
